$wb = $excel.ActiveWorkbook

# "line_imp" sheet (sheet2.xml) gets a new transformer column layout:
#   F1 "transformer" -> "t_x"
#   new G1 "t_a"
#   F2:F6 (previously "no") are cleared
$lineImp = $wb.Worksheets.Item("line_imp")
$lineImp.Range("F1").Value = "t_x"
$lineImp.Range("G1").Value = "t_a"
$lineImp.Range("F2:F6").ClearContents()

$initial = $wb.Worksheets.Item("initial")
$initial.Range("G5").Select()

$lineImp.Activate()
$lineImp.Range("F2").Select()
